$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for each row ---

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.817.27"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.21"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.72"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.53"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  +2.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.34"
$ws.Range("E12").Value = "  +3.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.398.72"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.780"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.078.41"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.811.67"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.05"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.87"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.78"
$ws.Range("E26").Value = "  +8.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.59"
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.65"
$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.33"
$ws.Range("E37").Value = "  -2.03%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("E39").Value = "  -0.99%  "

$ws.Range("E40").Value = "  +9.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.36"
$ws.Range("E41").Value = "  +2.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0971"
$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.19"
$ws.Range("E44").Value = "  +1.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.453.45"
$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.23"
$ws.Range("E49").Value = "  -3.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.281.46"
$ws.Range("E51").Value = "  +0.17%  "

# --- Swap rows 47 and 48: FTXToken <-> ARBITRUM ---
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.06"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.12"
$ws.Range("E48").Value = "  -4.91%  "

